$d = $word.ActiveDocument

# The "Totally Optional Reading & Viewing" list used to contain a plain-text,
# markdown-style pseudo-link:
#   [A Mockery of Justice for the Poor,] (http://www.nytimes.com/2016/04/30/opinion/a-mockery-of-justice-for-the-poor.html)
# Turn it into a real Word hyperlink, matching the style used by the rest of
# the list's hyperlinks (display text "A Mockery of Justice for the Poor,"
# with the "Hyperlink" character style, linking to the NYTimes article).

$target = "http://www.nytimes.com/2016/04/30/opinion/a-mockery-of-justice-for-the-poor.html"
$oldText = "[A Mockery of Justice for the Poor,] (" + $target + ")"
$displayText = "A Mockery of Justice for the Poor,"

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($rng, $target, "", "", $displayText) | Out-Null
}
